$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column B (shifts old B..G to D..I)
$ws.Columns("B:C").Insert()

# 2. Insert a new row before row 3 (adds a fresh blank row 3 below the existing data)
$ws.Rows("3:3").Insert()

# --- Row 1 (header) ---
$ws.Range("B1").Value = "Unnamed: 0.2"
$ws.Range("C1").Value = "Unnamed: 0.1"

# Copy the header style (from D1, which already carries the original header style) onto the
# two new header cells so B1/C1 match D1's bold/border/alignment formatting.
$ws.Range("D1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
# The column insert copied formatting from column A into B2/C2; clear it so they have no explicit style
$ws.Range("B2:C2").ClearFormats()

# --- Row 3 (new row) ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
# Empty-string text cells (not plain blanks): the leading "'" forces a text-typed, empty value.
$ws.Range("C3").Value = "'"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "Housing Purchase"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 100000
# Leading "'" keeps this as literal text instead of being parsed into a date serial number.
$ws.Range("I3").Value = "'1/12/2025"

# Copy style from A2 (s="1") onto A3 to match the rest of column A
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The rest of row 3 should carry no explicit style (the quote-prefix trick above can tag a
# style internally) - clear formats on everything except A3 to match the source file.
$ws.Range("B3:I3").ClearFormats()
